$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New repair-ticket rows for 2024-05-22, appended after the existing last row (168).
$newRows = @(
    @("2024-05-22", "11:48:35", "-", "-", "La cámara no detecta Top Cover", "-", "-", "11:48:43", "0:00:08"),
    @("2024-05-22", "11:48:47", "-", "-", "Detección de sealling mal puesto", "-", "-", "11:48:48", "0:00:01"),
    @("2024-05-22", "11:51:01", "-", "-", "-", "-", "Marco atascado en parte inferior", "11:51:05", "0:00:04"),
    @("2024-05-22", "11:51:07", "-", "-", "-", "-", "No detecta marcas Power", "11:51:10", "0:00:03"),
    @("2024-05-22", "11:53:26", "-", "Cámara no detecta Pcb", "-", "-", "-", "11:53:33", "0:00:07"),
    @("2024-05-22", "11:53:35", "-", "Cámara no detecta skeleton", "-", "-", "-", "11:53:39", "0:00:04"),
    @("2024-05-22", "11:54:07", "-", "Cámara no detecta skeleton", "-", "-", "-", "11:54:12", "0:00:05"),
    @("2024-05-22", "11:54:10", "-", "Cámara no detecta foam derecho", "-", "-", "-", "11:54:14", "0:00:04"),
    @("2024-05-22", "11:54:51", "-", "Cámara no detecta foams", "-", "-", "-", "11:54:53", "0:00:02"),
    @("2024-05-22", "11:54:56", "-", "Etiquetadora", "-", "-", "-", "11:55:00", "0:00:04"),
    @("2024-05-22", "11:55:35", "-", "Power atascado en prensa, cuesta sacar", "-", "-", "-", "11:55:37", "0:00:02"),
    @("2024-05-22", "11:56:09", "-", "Cámara no detecta foam derecho", "-", "-", "-", "11:56:10", "0:00:01"),
    @("2024-05-22", "11:56:15", "-", "AOI (malla)", "-", "-", "-", "11:56:16", "0:00:01")
)

$startRow = 169

# Column A holds date-looking text ("2024-05-22"); force it to text so Excel
# doesn't silently convert it to a date serial number.
$endRow = $startRow + $newRows.Count - 1
$ws.Range("A" + $startRow + ":A" + $endRow).NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowData = $newRows[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}
